$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 8 (i2o-dev-30watt) execution status from "No" to "Yes"
$ws.Range("B8").Value = "Yes"

# Add new client rows
$ws.Range("A9").Value = "i2o-dev-ossur"
$ws.Range("B9").Value = "No"

$ws.Range("A10").Value = "i2o-preprod-mycharge"
$ws.Range("B10").Value = "No"

$ws.Range("A11").Value = "i2o-preprod-interaxon"
$ws.Range("B11").Value = "No"

$ws.Range("A12").Value = "i2o-dev-ffl"
$ws.Range("B12").Value = "No"

$ws.Range("A13").Value = "i2o-preprod-talkworks"
$ws.Range("B13").Value = "No"

# Copy style from existing data row to new rows
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B13").PasteSpecial(-4122) | Out-Null

# Update selection to match final state
$ws.Range("B9").Select()
